# "Add files via upload" — bmi.xlsx re-upload that cleared the helper
# BMI-category column (H) while leaving everything else (labels, weights,
# heights, styles) untouched, and left the selection sitting on H16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously-computed BMI category text out of H2:H21. The
# cell's style (s="1"/s="2") stays in place — only the string value
# (and its "t=s" type marker) goes away. Saving afterwards naturally
# compacts the shared-string table, which drops the now-unused numeric
# height strings and shifts "U" down to the first free slot.
$ws.Range("H2:H21").ClearContents()

# Leave the selection where the user ended up after clearing the column.
$ws.Range("H16").Select()
